$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.428.82'
$ws.Range("E2").Value = '  +1.71%  '

$ws.Range("D3").Value = '2.017.18'
$ws.Range("E3").Value = '  +5.96%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5128'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4225'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08726'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.97%  '

$ws.Range("E11").Value = '  +3.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.98%  '

$ws.Range("D13").Value = '2.019.37'
$ws.Range("E13").Value = '  +5.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.610'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.491'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.46'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001116'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06534'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.72%  '

$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.223'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.14%  '

$ws.Range("D23").Value = '30.491.15'
$ws.Range("E23").Value = '  +1.80%  '

$ws.Range("E24").Value = '  +6.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.229'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.53%  '

$ws.Range("E26").Value = '  +5.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.425'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.141'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.98%  '

$ws.Range("E32").Value = '  +2.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.084'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.832'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.06%  '

$ws.Range("E35").Value = '  +15.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02539'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.487'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06669'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.28'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2205'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.064'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6679'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.234'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.86%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6195'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.193'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.660'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.82%  '

$ws.Range("E49").Value = '  +5.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.76%  '
